$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.177.52"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.835.93"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'241.92"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "'0.6606"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D8").Value = "'0.07421"
$ws.Range("E8").Value = "  +1.33%  "
$ws.Range("D9").Value = "'0.2936"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").Value = "'0.07757"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").Value = "1.839.16"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "'4.989"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "'0.6668"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "'82.93"
$ws.Range("E15").Value = "  -3.51%  "
$ws.Range("D16").Value = "'6.106"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "'0.000008546"
$ws.Range("E17").Value = "  +4.30%  "
$ws.Range("D18").Value = "29.177.60"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").Value = "2.107.85"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").Value = "'227.06"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'7.098"
$ws.Range("E23").Value = "  -2.12%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'159.49"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("D26").Value = "'8.613"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").Value = "'0.1399"
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("D28").Value = "'17.96"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("D30").Value = "'4.114"
$ws.Range("E30").Value = "  -2.54%  "
$ws.Range("D31").Value = "'4.042"
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D33").Value = "'0.05269"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").Value = "'1.860"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").Value = "'0.7360"
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("D36").Value = "'1.146"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("D37").Value = "'2.659"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").Value = "1.298.93"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").Value = "'0.01793"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("D41").Value = "'0.9197"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "'0.08543"
$ws.Range("E43").Value = "  +17.99%  "
$ws.Range("D44").Value = "'0.9994"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "'102.87"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").Value = "1.994.92"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").Value = "'0.5140"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").Value = "'63.55"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'0.05844"
$ws.Range("E51").Value = "  -1.11%  "
